$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column G ("Recorded By") values need the order of the comma-separated
# recorder names reversed for the rows below (data refresh / re-sort of
# recorder attribution during daily attendance processing).
$values = @{
    2 = 'backup@backdoor.com, System, system'
    3 = 'System, dnasr281@gmail.com'
    4 = 'backup@backdoor.com, System'
    5 = 'backup@backdoor.com, System'
    6 = 'System, dnasr281@gmail.com'
    10 = 'System, dnasr281@gmail.com'
    11 = 'System, dnasr281@gmail.com'
    12 = 'System, dnasr281@gmail.com'
    13 = 'System, dnasr281@gmail.com'
    14 = 'System, dnasr281@gmail.com'
    15 = 'System, dnasr281@gmail.com'
    17 = 'System, dnasr281@gmail.com'
    18 = 'System, dnasr281@gmail.com'
    29 = 'backup@backdoor.com, System, system'
    30 = 'System, dnasr281@gmail.com'
    31 = 'backup@backdoor.com, System'
    32 = 'backup@backdoor.com, System'
    33 = 'System, dnasr281@gmail.com'
    37 = 'System, dnasr281@gmail.com'
    38 = 'System, dnasr281@gmail.com'
    39 = 'System, dnasr281@gmail.com'
    40 = 'System, dnasr281@gmail.com'
    41 = 'System, dnasr281@gmail.com'
    42 = 'System, dnasr281@gmail.com'
    44 = 'System, dnasr281@gmail.com'
    45 = 'System, dnasr281@gmail.com'
    56 = 'backup@backdoor.com, System, system'
    57 = 'System, dnasr281@gmail.com'
    58 = 'backup@backdoor.com, System'
    59 = 'backup@backdoor.com, System'
    60 = 'System, dnasr281@gmail.com'
    64 = 'System, dnasr281@gmail.com'
    65 = 'System, dnasr281@gmail.com'
    66 = 'System, dnasr281@gmail.com'
    67 = 'System, dnasr281@gmail.com'
    68 = 'System, dnasr281@gmail.com'
    69 = 'System, dnasr281@gmail.com'
    71 = 'System, dnasr281@gmail.com'
    72 = 'System, dnasr281@gmail.com'
    84 = 'backup@backdoor.com, System'
    85 = 'backup@backdoor.com, System'
    86 = 'System, dnasr281@gmail.com'
    87 = 'System, dnasr281@gmail.com'
    88 = 'System, dnasr281@gmail.com'
    89 = 'System, dnasr281@gmail.com'
    93 = 'System, dnasr281@gmail.com'
    95 = 'System, dnasr281@gmail.com'
    96 = 'System, dnasr281@gmail.com'
    110 = 'backup@backdoor.com, System'
    111 = 'backup@backdoor.com, System'
    112 = 'System, dnasr281@gmail.com'
    113 = 'System, dnasr281@gmail.com'
    114 = 'System, dnasr281@gmail.com'
    115 = 'System, dnasr281@gmail.com'
    119 = 'System, dnasr281@gmail.com'
    121 = 'System, dnasr281@gmail.com'
    122 = 'System, dnasr281@gmail.com'
    136 = 'backup@backdoor.com, System'
    137 = 'backup@backdoor.com, System'
    138 = 'System, dnasr281@gmail.com'
    139 = 'System, dnasr281@gmail.com'
    140 = 'System, dnasr281@gmail.com'
    141 = 'System, dnasr281@gmail.com'
    145 = 'System, dnasr281@gmail.com'
    147 = 'System, dnasr281@gmail.com'
    148 = 'System, dnasr281@gmail.com'
}

foreach ($row in $values.Keys) {
    $ws.Cells.Item($row, 7).Value = $values[$row]
}
